# 19,12,2024 dodane dwa zestawy danych, aor1_rna do poprawy
#
# Adds two new dataset rows (abrc3, absaA2_scoe) and one more (aor1_rna)
# into the existing "data database" sheet, filling in the "data", "type of
# experiment" and "specifics of the experiment" columns for rows 8-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (record #6): abrC3 dataset
$ws.Range("B8").Value = "abrc3"
$ws.Range("C8").Value = "RNA-seq"
$ws.Range("D8").Value = "abrC3 mutant relative to the parent strain M145"

# Row 9 (record #7): absaA2_scoe dataset (type of experiment left blank, as in source)
$ws.Range("B9").Value = "absaA2_scoe"
$ws.Range("D9").Value = "troche nie wiem co z tym zrobić"

# Row 10 (record #8): aor1_rna dataset
$ws.Range("B10").Value = "aor1_rna"
$ws.Range("C10").Value = "RNA-seq"
$ws.Range("D10").Value = "Genes whose transcript levels changed at least 3-fold (FC) in S. coelicolor ∆aor1 compared with those of M145"

# D8/D10 lose their fill+border ("fill" alignment only, like D6/D7 had);
# D9 keeps the bordered/filled look. Copy the already-plain style from D7.
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the author's final selection/cursor position.
$ws.Range("G9").Select()
